$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add four new header columns for the trip-upload template:
# Actor Shipper, Container Number, Shipper Reference, Booking number
$ws.Range("J1").Value = "Actor Shipper"
$ws.Range("K1").Value = "Container Number"
$ws.Range("L1").Value = "Shipper Reference"
$ws.Range("M1").Value = "Booking number"

# Size the new columns to match the published template
$ws.Columns.Item(10).ColumnWidth = 9.830729166666666
$ws.Columns.Item(11).ColumnWidth = 14.166666666666666
$ws.Columns.Item(12).ColumnWidth = 12.330729166666666
$ws.Columns.Item(13).ColumnWidth = 13.166666666666666

# Shipper Reference / Booking number headers use the plain (non-wrapping)
# bordered header style, matching the rest of the template
$plainHeader = $ws.Range("L1:M1")
$plainHeader.Font.Name = "Open Sans"
$plainHeader.Font.Size = 8
$plainHeader.HorizontalAlignment = -4108
$plainHeader.VerticalAlignment = -4108
$plainHeader.WrapText = $false
$plainHeader.Borders.LineStyle = 1
$plainHeader.Borders.Color = 0

# Reflect the author's view state after adding the new columns
$ws.Range("K5").Select()
$excel.ActiveWindow.Zoom = 254
